$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "61.953.05"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -0.88%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.410.02"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -0.66%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  +0.12%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "409.71"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +0.51%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "129.07"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("E7").Value = "  +6.17%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +5.08%  "

$ws.Range("E10").Value = "  -1.65%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "43.13"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +1.86%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000220"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +28.59%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.33"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +9.54%  "

$ws.Range("E14").Value = "  -0.34%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.39"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +7.42%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.950.63"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -0.93%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.406.29"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -0.30%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.52"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +8.40%  "

$ws.Range("E19").Value = "  +6.63%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "61.969.81"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.85%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "448.93"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +42.34%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "91.70"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +8.19%  "

$ws.Range("E23").Value = "  +0.26%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.21"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +2.42%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.30"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +3.59%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.46"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +15.50%  "

$ws.Range("E27").Value = "  +10.97%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.79"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +1.11%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.70"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -1.41%  "

$ws.Range("E30").Value = "  -0.26%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.02"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +4.96%  "

$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("E33").Value = "  -0.97%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.72"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -3.72%  "

$ws.Range("E35").Value = "  +0.00%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0507"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +4.00%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "53.83"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +3.96%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("E40").Value = "  +8.05%  "

$ws.Range("E41").Value = "  -0.48%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.320"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -1.02%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "143.90"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("E44").Value = "  +9.95%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.01"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("E46").Value = "  +14.45%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "16.63"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -2.30%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.154"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +25.95%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.57"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +5.52%  "

$ws.Range("E50").Value = "  +7.67%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.751.49"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -0.77%  "
